# Applies the README edits described by the commit diff:
#  - Collapses a handful of runs that were only split apart by inline
#    spell/grammar-check markers (<w:proofErr/>) back into single runs,
#    and drops the now-stale proofErr markers themselves.
#  - Replaces the placeholder "TODO!!!!!!!" paragraph (Output section)
#    with the real output-file note, split across two paragraphs (one
#    plain sentence, one bold bulleted NOTE).

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaByText {
    param([string]$needle)
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    throw "Paragraph containing '$needle' not found"
}

# 1) "... in the command line (cmd):" - merge " (" + "cmd" + ")" into
#    a single run and drop the spellStart/spellEnd proofErr pair.
$p1 = Get-ParaByText "in the command line"
$xml1 = '<w:p ' + $wordNs + '><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Before running the script, ensure you have the necessary packages installed. You can install them using the following commands</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> in the command line</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (cmd)</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# 2) "pip install pandas" - drop the gramStart/gramEnd proofErr pair.
$p2 = Get-ParaByText "pip install "
$p2 = Get-ParaByText "pandas"
$xml2 = '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">pip install </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>pandas</w:t></w:r></w:p>'
$p2.Range.InsertXML($xml2)

# 3) "pip install openpyxl" - drop the spellStart/gramStart/spellEnd/gramEnd
#    proofErr markers.
$p3 = Get-ParaByText "openpyxl"
$xml3 = '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">pip install </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>openpyxl</w:t></w:r></w:p>'
$p3.Range.InsertXML($xml3)

# 4) "(Google chrome (the standard webbrowser) needs to be downloaded for
#    chromedriver to function)" - merge the 5 runs split by spellStart/
#    spellEnd proofErr pairs back into a single run.
$p4 = Get-ParaByText "Google chrome"
$xml4 = '<w:p ' + $wordNs + '>' + `
    '<w:r><w:t xml:space="preserve">The script uses the Chrome WebDriver. </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Make sure you have Google Chrome installed on your machine</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' + `
    '<w:r><w:t>(Google chrome (the standard webbrowser) needs to be downloaded for chromedriver to function)</w:t></w:r></w:p>'
$p4.Range.InsertXML($xml4)

# 5) "Two backslashes (\\) ... is actually just a character." - drop the
#    gramStart/gramEnd proofErr pair and merge the "  \" run with the
#    " is just defining that the 2" run.
$p5 = Get-ParaByText "Two backslashes"
$xml5 = '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Two backslashes (\\) are required for to define a backslash character (\) in a python string. This is because \ is an escape character in python strings (meaning it allows us to note special character, such as \n newline character). Therefore the 1</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  \ is just defining that the 2</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> \ is actually just a character.</w:t></w:r></w:p>'
$p5.Range.InsertXML($xml5)

# 6) Replace the "TODO!!!!!!!" placeholder paragraph with the real output
#    note, split into two paragraphs: a plain sentence describing the csv
#    output, and a bold bulleted NOTE warning not to have the file open.
$lq = [char]0x2018
$rq = [char]0x2019
$outputsText = "Outputs to a csv called " + $lq + "scraped_auto_premium.csv" + $rq + "."

$p6 = Get-ParaByText "TODO"
$xml6 = '<w:p ' + $wordNs + '><w:r><w:t>' + $outputsText + '</w:t></w:r></w:p>' + `
    '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>NOTE: MAKE SURE this file is NOT OPEN while the code is running</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> or the scraped premiums cannot be written to it.</w:t></w:r></w:p>'
$p6.Range.InsertXML($xml6)

Write-Host "Edits applied"
